$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Wipe out the old data rows (2-3) completely so leftover formatting
#      (row height, stray Helvetica Neue font on E2, etc.) is gone, and new
#      rows start from a clean default style ----
$ws.Range("A2:E3").Delete()

# ---- Header row (row 1) : keep values, fix formatting below ----
$ws.Cells.Item(1,1).Value = "Vessel Type"
$ws.Cells.Item(1,2).Value = "Barcode"
$ws.Cells.Item(1,3).Value = "UMI Length"
$ws.Cells.Item(1,4).Value = "Spacer Length"
$ws.Cells.Item(1,5).Value = "Location"

# ---- Data rows 2-8 (fresh cells, default style) ----
$ws.Cells.Item(2,1).Value = "Eppendorf96"
$ws.Cells.Item(2,2).Value = 12345
$ws.Cells.Item(2,3).Value = 6
$ws.Cells.Item(2,4).Value = 3
$ws.Cells.Item(2,5).Value = "Inline First Read"

$ws.Cells.Item(3,1).Value = "Eppendorf96"
$ws.Cells.Item(3,2).Value = 34567
$ws.Cells.Item(3,3).Value = 3
$ws.Cells.Item(3,4).Value = 1
$ws.Cells.Item(3,5).Value = "Before Second Index Read"

$ws.Cells.Item(4,1).Value = "Eppendorf96"
$ws.Cells.Item(4,2).Value = 66789
$ws.Cells.Item(4,3).Value = 9
$ws.Cells.Item(4,4).Value = 1
$ws.Cells.Item(4,5).Value = "Inline Second Read"

$ws.Cells.Item(5,1).Value = "Eppendorf96"
$ws.Cells.Item(5,2).Value = 77891
$ws.Cells.Item(5,3).Value = 9
$ws.Cells.Item(5,4).Value = 3
$ws.Cells.Item(5,5).Value = "Inline Second Read"

$ws.Cells.Item(6,1).Value = "MatrixTube075"
$ws.Cells.Item(6,2).Value = 87654
$ws.Cells.Item(6,3).Value = 3
$ws.Cells.Item(6,4).Value = 2
$ws.Cells.Item(6,5).Value = "Before First Read"

$ws.Cells.Item(7,1).Value = "MatrixTube075"
$ws.Cells.Item(7,2).Value = 87654
$ws.Cells.Item(7,3).Value = 3
$ws.Cells.Item(7,4).Value = 2
$ws.Cells.Item(7,5).Value = "Before Second Read"

$ws.Cells.Item(8,1).Value = "Eppendorf96"
$ws.Cells.Item(8,2).Value = 77891
$ws.Cells.Item(8,3).Value = 2
$ws.Cells.Item(8,4).Value = 2
$ws.Cells.Item(8,5).Value = "Before First Read"

# ---- Header formatting ----
# A1: non-bold Arial 10, color #222222
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Color = 2236962
$ws.Range("A1").Font.Bold = $false

# B1:E1: bold Arial 10, color #222222
$ws.Range("B1:E1").Font.Name = "Arial"
$ws.Range("B1:E1").Font.Size = 10
$ws.Range("B1:E1").Font.Color = 2236962
$ws.Range("B1:E1").Font.Bold = $true

# ---- Row 3 C:D bold ----
$ws.Range("C3:D3").Font.Bold = $true

# ---- Selection matches saved file state ----
$ws.Range("E8").Select()
